# Daily refresh of the cryptos price table (GitHub Actions job).
# Updates Price (col D) / Volume(1h) (col E) for each coin row, and
# swaps the Hedera / VeChain rows (38 <-> 39) which changed rank order.
# Numeric-looking price strings are written with a leading "'" so Excel
# keeps them as text (matching the sheet's original inlineStr/text cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.873.05"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "3.146.42"
$ws.Range("E3").Value = "  +1.99%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'572.02"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").Value = "'151.37"
$ws.Range("E6").Value = "  +4.60%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "3.145.61"
$ws.Range("E8").Value = "  +2.02%  "

$ws.Range("E9").Value = "  +4.30%  "

$ws.Range("E10").Value = "  +5.10%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").Value = "'0.502"
$ws.Range("E12").Value = "  +6.63%  "

$ws.Range("E13").Value = "  +11.18%  "

$ws.Range("D14").Value = "'37.50"
$ws.Range("E14").Value = "  +6.74%  "

$ws.Range("D15").Value = "3.657.63"
$ws.Range("E15").Value = "  +2.02%  "

$ws.Range("D16").Value = "64.929.90"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").Value = "'7.22"
$ws.Range("E17").Value = "  +6.89%  "

$ws.Range("D18").Value = "3.144.01"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").Value = "'512.19"
$ws.Range("E20").Value = "  +6.80%  "

$ws.Range("D21").Value = "'14.93"
$ws.Range("E21").Value = "  +7.19%  "

$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "  +9.49%  "

$ws.Range("D23").Value = "'15.31"
$ws.Range("E23").Value = "  +10.55%  "

$ws.Range("D24").Value = "'7.82"
$ws.Range("E24").Value = "  +3.89%  "

$ws.Range("D25").Value = "'84.98"
$ws.Range("E25").Value = "  +4.74%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  +3.88%  "

$ws.Range("E28").Value = "  +8.58%  "

$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +5.67%  "

$ws.Range("D30").Value = "'27.97"
$ws.Range("E30").Value = "  +6.69%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  +3.54%  "

$ws.Range("E33").Value = "  +6.43%  "

$ws.Range("D34").Value = "'6.08"
$ws.Range("E34").Value = "  +8.62%  "

$ws.Range("E35").Value = "  +6.06%  "

$ws.Range("D36").Value = "'55.41"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").Value = "'482.50"
$ws.Range("E37").Value = "  +5.47%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.0862"
$ws.Range("E38").Value = "  +4.82%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0423"
$ws.Range("E39").Value = "  +3.54%  "

$ws.Range("D40").Value = "'3.01"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "3.118.09"
$ws.Range("E41").Value = "  +4.33%  "

$ws.Range("E42").Value = "  +4.69%  "

$ws.Range("E43").Value = "  +4.20%  "

$ws.Range("D44").Value = "'0.291"
$ws.Range("E44").Value = "  +11.27%  "

$ws.Range("E45").Value = "  +14.41%  "

$ws.Range("D46").Value = "'29.21"
$ws.Range("E46").Value = "  +4.13%  "

$ws.Range("D47").Value = "0.0₃0572"
$ws.Range("E47").Value = "  +10.88%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("E50").Value = "  +10.51%  "

$ws.Range("D51").Value = "'118.92"
$ws.Range("E51").Value = "  -1.64%  "
